$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: the reporting period start date moved forward (B8) ---
$ws.Range("B8").Value = 44105

# --- L8: replace the "awaiting financial close" note with the final note ---
$ws.Range("L8").Value = "Los campos en blanco, obedecen a que en este Organismo, no se generaron bajas de bienes muebles durante el periodo que se reporta. "
$ws.Range("L8").HorizontalAlignment = -4131   # xlLeft (was justify)
$ws.Range("L8").WrapText = $true

# Row 8 is shorter now that the note text is shorter
$ws.Rows("8").RowHeight = 45

# Row 3's custom height is no longer needed
[void]$ws.Rows("3").AutoFit()

# --- New row 9: a stray bold red cell left behind at H9 ---
$ws.Range("H9").Font.Bold = $true
$ws.Range("H9").Font.Color = 255

# Leave the selection where the author left off
[void]$ws.Range("D14").Select()

Write-Host "done"
